$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.491.93"
$ws.Range("E2").Value = "  -6.89%  "
$ws.Range("D3").Value = "2.887.06"
$ws.Range("E3").Value = "  -5.31%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.08"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.20%  "
$ws.Range("D8").Value = "2.878.12"
$ws.Range("E8").Value = "  -5.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("E10").Value = "  -10.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.73"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -10.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.433"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000210"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -10.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.119"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "3.363.81"
$ws.Range("E16").Value = "  -5.17%  "
$ws.Range("D17").Value = "2.878.86"
$ws.Range("E17").Value = "  -5.38%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.49"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "57.402.55"
$ws.Range("E19").Value = "  -7.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "407.66"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.20%  "
$ws.Range("E21").Value = "  -5.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.651"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.70"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -9.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "76.60"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.53%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.46"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.89"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.41%  "
$ws.Range("E32").Value = "  -5.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0936"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.01"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -13.72%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.10%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.889"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.84%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.24"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.28"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.31%  "
$ws.Range("D39").Value = "0.0₃0613"
$ws.Range("E39").Value = "  -12.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0343"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.22%  "
$ws.Range("E41").Value = "  -4.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "361.47"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.84%  "
$ws.Range("D43").Value = "2.598.06"
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("E45").Value = "  -8.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "118.24"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.20%  "
$ws.Range("E47").Value = "  -5.14%  "
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.91"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.09"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.28%  "
